$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.179.08'
$ws.Range('E2').Value = '  +0.96%  '
$ws.Range('D3').Value = '3.632.97'
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '195.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '574.10'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('D7').Value = '3.628.58'
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').Value = '  +1.51%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('E11').Value = '  +4.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.56%  '
$ws.Range('E13').Value = '  +15.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.12'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.11%  '
$ws.Range('D15').Value = '4.206.26'
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').Value = '3.638.80'
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.55'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('D19').Value = '68.110.43'
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '402.62'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +22.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.01'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.84%  '
$ws.Range('E26').Value = '  +3.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.20%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.16%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.85'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.12'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +19.53%  '
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.67'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '691.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +17.67%  '
$ws.Range('E34').Value = '  +2.20%  '
$ws.Range('E35').Value = '  +5.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.81'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.68'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.00%  '
$ws.Range('E38').Value = '  +10.22%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '0.0₃0796'
$ws.Range('E40').Value = '  +7.76%  '
$ws.Range('E41').Value = '  +16.81%  '
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.12'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +11.75%  '
$ws.Range('D44').Value = '3.178.39'
$ws.Range('E44').Value = '  +18.21%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('E46').Value = '  +24.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0422'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('E48').Value = '  +1.99%  '
$ws.Range('E49').Value = '  +5.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '142.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.24%  '
